$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(19).Insert()

Write-Output "Inserted"
Write-Output $ws.Cells.Item(19, 1).Value()
Write-Output $ws.Cells.Item(20, 3).Value()
Write-Output $ws.Cells.Item(34, 16).Value()
Write-Output $ws.Cells.Item(35, 16).Value()
